$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 and B3 keep their existing text values but pick up the bordered
# style that A2:A4 / C2:C4 already use (xf index 1).
$ws.Range("A2").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)   # xlPasteFormats

# B4's value changes to a new job number. Force it to be stored as text
# (matching B2/B3) with a leading apostrophe, then re-apply the bordered
# format so the cell ends up using the same style as B2/B3 instead of a
# brand-new "quote prefixed" style.
$ws.Range("B4").Value = "'32339569"
$ws.Range("A2").Copy()
$ws.Range("B4").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0

# The trailing blank row (row 5) is removed entirely.
$ws.Rows.Item(5).Delete()

# Mirror the author's final selection.
$ws.Range("E12").Select()
